# Update the 'time_taken' timestamps on the existing 'data' sheet to
# reflect the later re-run time captured in the diff.
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$timeTaken = @(
    "2021-10-05 14:20:28.341250",
    "2021-10-05 14:20:28.341257",
    "2021-10-05 14:20:28.341261",
    "2021-10-05 14:20:28.341263",
    "2021-10-05 14:20:28.341266",
    "2021-10-05 14:20:28.341268",
    "2021-10-05 14:20:28.341271",
    "2021-10-05 14:20:28.341273",
    "2021-10-05 14:20:28.341276",
    "2021-10-05 14:20:28.341279",
    "2021-10-05 14:20:28.341281",
    "2021-10-05 14:20:28.341284"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# Add a new 'metadata' worksheet right after the 'data' sheet describing
# the panel query that produced this workbook.
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Reuse the bordered/bold header formatting already defined in the
# workbook (as seen on the 'data' sheet's header row / index column)
# instead of building a brand-new style.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Generalised pustular psoriasis"
$meta.Cells.Item(2, 3).Value = 132

# Keep the version string as literal text ("1.9"), not a coerced number.
$versionCell = $meta.Cells.Item(2, 4)
$versionCell.NumberFormat = "@"
$versionCell.Value = "1.9"
$versionCell.ClearFormats()

$meta.Cells.Item(2, 5).Value = "2021-01-29T15:24:57.026812Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:20:28.337684"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/132/?format=json"

# Leave the originally active 'data' sheet selected, matching the
# workbook's prior view state.
$data.Activate()
